# Revert "Merge pull request #493 ..." — re-adds the "Penalty Issued" column
# (column I) to the Dairy Test Threshold Template report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: width to roughly match the restored column ---
$ws.Columns.Item(9).ColumnWidth = 10.5

# --- I4: new header cell "Penalty Issued", formatted like the other
#         row-4 header cells (bold, centered, wrapped) ---
$ws.Range("I4").Value = "Penalty Issued"
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- I6: template placeholder text for the penalty-issued value ---
$ws.Range("I6").Value = "{d.Reg[i].PenaltyIssued}"

# --- I3 / I5: blank spacer cells carrying the same bottom border used
#              elsewhere in the template ---
foreach ($addr in @("I3", "I5")) {
    $bottom = $ws.Range($addr).Borders.Item(9)   # xlEdgeBottom
    $bottom.Color = 0                            # black
    $bottom.Weight = 2                           # xlThin
    $bottom.LineStyle = 1                        # xlContinuous
}

# --- Selection moved as part of the resaved view state ---
$ws.Range("G16").Select()
